$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.StandardWidth = 11.19921875
